$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.455.44"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "3.172.25"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'238.23"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "'624.11"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +5.68%  "
$ws.Range("D8").Value = "'0.370"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "3.174.60"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("D11").Value = "'0.742"
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'35.42"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "90.819.62"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "3.764.96"
$ws.Range("D18").Value = "3.185.73"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").Value = "'3.71"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "'15.10"
$ws.Range("E20").Value = "  +9.10%  "
$ws.Range("D21").Value = "'5.89"
$ws.Range("E21").Value = "  +7.61%  "
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").Value = "'441.89"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").Value = "'5.76"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").Value = "'89.37"
$ws.Range("E26").Value = "  +10.33%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("D28").Value = "3.329.88"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'0.129"
$ws.Range("E30").Value = "  +48.48%  "
$ws.Range("D31").Value = "'0.232"
$ws.Range("E31").Value = "  +20.03%  "
$ws.Range("D32").Value = "'0.170"
$ws.Range("E32").Value = "  +8.12%  "
$ws.Range("D33").Value = "'9.48"
$ws.Range("E33").Value = "  +3.69%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("E35").Value = "  +11.55%  "
$ws.Range("D36").Value = "'7.79"
$ws.Range("E36").Value = "  +10.57%  "
$ws.Range("D37").Value = "'26.45"
$ws.Range("D38").Value = "'506.60"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("E40").Value = "  +8.19%  "
$ws.Range("D41").Value = "'0.452"
$ws.Range("E41").Value = "  +12.77%  "
$ws.Range("D42").Value = "'3.76"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").Value = "'3.42"
$ws.Range("E43").Value = "  -9.09%  "
$ws.Range("D44").Value = "'22.12"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D46").Value = "'0.725"
$ws.Range("E46").Value = "  +7.40%  "
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "'1.92"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E49").Value = "  +5.17%  "
$ws.Range("D50").Value = "'4.43"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "'44.03"
$ws.Range("E51").Value = "  -0.96%  "
